$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("B8").Value = 16
